# "full version of v4, still a little slower T_T"
#
# The existing "remark / git-version" column lived in column N, aligned with
# the FLK/base_channel tables that occupy columns B:L. A new numeric column
# of data is being introduced right where that remark column used to be, so
# the remark column itself is pushed one slot to the right (N -> O).
#
# A brand-new benchmark row ("加入v4") is inserted right after the existing
# 256-channel FLKv3 row (old row 17), which also gains four new measurement
# values (J:M). Everything below (the base_channel/model/amount table that
# used to start at row 19) shifts down by one row to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert a new column before N, and a new row before 19.
$ws.Columns("N").Insert()
$ws.Rows("19").Insert()

# 2) The FLKv3/256-channel row (17) picks up four new benchmark numbers.
$ws.Range("J17").Value = 0.029
$ws.Range("K17").Value = 0.255
$ws.Range("L17").Value = 0.174
$ws.Range("M17").Value = 0.159

# 3) Fill in the brand-new row 18 ("加入v4").
$ws.Range("A18").Value = "加入v4"
$ws.Range("C18").Value = 256
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 128
$ws.Range("F18").Value = 112
$ws.Range("G18").Value = 112
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 100
$ws.Range("J18").Value = 0.029
$ws.Range("K18").Value = 0.253
$ws.Range("L18").Value = 0.179
$ws.Range("M18").Value = 0.169
$ws.Range("N18").Value = 0.035
$ws.Range("O18").Value = "22010a2"

# 4) Keep the selection where the author left it.
$ws.Range("O18").Select() | Out-Null
